$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeLoginCredentials")

# Correct the generated usernames (Username column, C2:C4) -- append the
# missing "2" so each login lines up with the fixed test data.
$ws.Range("C2").Value = "barhn45612"
$ws.Range("C3").Value = "mayrw45612"
$ws.Range("C4").Value = "Akbr45612"

# Move the sheet's saved selection from the stale C8 (outside the used
# range) back onto the table at C4.
$ws.Activate()
$ws.Range("C4").Select()
